$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.273.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.663.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.79%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5329'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2634'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06357'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07817'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.891.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5529'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8193'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.678'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.029'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.183'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05875'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.278'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.587'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.609'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.821'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.423'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5795'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01603'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8634'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.832'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.009'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.046.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.801.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.013'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4379'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.020'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05162'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.429'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.37%  '
